# Scheduled runner update: refresh market-price-derived figures (currentAveragePrice*,
# LevePriceNQ/HQ, LeveProfitNQ/HQ) across the per-class Sheets, cell by cell.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 659.2
$ws.Range("J49").Value = 998.5
$ws.Range("L49").Value = 2995.5
$ws.Range("N49").Value = -3267.5
$ws.Range("H62").Value = 4891.727
$ws.Range("I62").Value = 4589.647
$ws.Range("J62").Value = 5918.8
$ws.Range("K62").Value = 4589.647
$ws.Range("L62").Value = 5918.8
$ws.Range("M62").Value = -3965.647
$ws.Range("N62").Value = -7166.8
$ws.Range("H65").Value = 4891.727
$ws.Range("I65").Value = 4589.647
$ws.Range("J65").Value = 5918.8
$ws.Range("K65").Value = 22948.235
$ws.Range("L65").Value = 29594
$ws.Range("M65").Value = -19828.235
$ws.Range("N65").Value = -35834
$ws.Range("H69").Value = 563804.75
$ws.Range("I69").Value = 778460.4399999999
$ws.Range("J69").Value = 5699.8
$ws.Range("K69").Value = 2335381.32
$ws.Range("L69").Value = 17099.4
$ws.Range("M69").Value = -2334507.32
$ws.Range("N69").Value = -18847.4
$ws.Range("H72").Value = 563804.75
$ws.Range("I72").Value = 778460.4399999999
$ws.Range("J72").Value = 5699.8
$ws.Range("K72").Value = 7006143.959999999
$ws.Range("L72").Value = 51298.2
$ws.Range("M72").Value = -7001775.959999999
$ws.Range("N72").Value = -60034.2
$ws.Range("H74").Value = 12389.421
$ws.Range("I74").Value = 2749.75
$ws.Range("K74").Value = 2749.75
$ws.Range("M74").Value = -1813.75
$ws.Range("H77").Value = 12389.421
$ws.Range("I77").Value = 2749.75
$ws.Range("K77").Value = 13748.75
$ws.Range("M77").Value = -9068.75
$ws.Range("H92").Value = 44707.156
$ws.Range("I92").Value = 18234.5
$ws.Range("J92").Value = 230015.75
$ws.Range("K92").Value = 18234.5
$ws.Range("L92").Value = 230015.75
$ws.Range("M92").Value = -16986.5
$ws.Range("N92").Value = -232511.75
$ws.Range("H104").Value = 138.4
$ws.Range("I104").Value = 138.4
$ws.Range("K104").Value = 415.2
$ws.Range("M104").Value = 1331.8
$ws.Range("H111").Value = 1231.8182
$ws.Range("I111").Value = 865
$ws.Range("K111").Value = 2595
$ws.Range("M111").Value = 472
$ws.Range("H113").Value = 2550
$ws.Range("I113").Value = 2400
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2400
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 854
$ws.Range("N113").Value = -9508
$ws.Range("H116").Value = 6656.4546
$ws.Range("I116").Value = 5778.7896
$ws.Range("J116").Value = 7847.5713
$ws.Range("K116").Value = 5778.7896
$ws.Range("L116").Value = 7847.5713
$ws.Range("M116").Value = -2336.7896
$ws.Range("N116").Value = -14731.5713
$ws.Range("H132").Value = 1781.0186
$ws.Range("I132").Value = 1669.2653
$ws.Range("K132").Value = 5007.7959
$ws.Range("M132").Value = -2477.7959
$ws.Range("H138").Value = 3454.9854
$ws.Range("I138").Value = 1511.25
$ws.Range("K138").Value = 4533.75
$ws.Range("M138").Value = 606.25
$ws.Range("H141").Value = 1712.8334
$ws.Range("I141").Value = 1675.4
$ws.Range("J141").Value = 1900
$ws.Range("K141").Value = 5026.200000000001
$ws.Range("L141").Value = 5700
$ws.Range("M141").Value = 153.7999999999993
$ws.Range("N141").Value = -16060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 6794
$ws.Range("I19").Value = 1028.5
$ws.Range("J19").Value = 18325
$ws.Range("K19").Value = 1028.5
$ws.Range("L19").Value = 18325
$ws.Range("M19").Value = -799.5
$ws.Range("N19").Value = -18783
$ws.Range("H32").Value = 17045.89
$ws.Range("I32").Value = 17089.012
$ws.Range("K32").Value = 17089.012
$ws.Range("M32").Value = -16802.012
$ws.Range("H45").Value = 645.4
$ws.Range("I45").Value = 604.0222
$ws.Range("K45").Value = 604.0222
$ws.Range("M45").Value = -227.0222
$ws.Range("H61").Value = 4278.3584
$ws.Range("I61").Value = 2818.8684
$ws.Range("K61").Value = 2818.8684
$ws.Range("M61").Value = -2606.8684
$ws.Range("H63").Value = 3019.8845
$ws.Range("I63").Value = 2333.6191
$ws.Range("K63").Value = 2333.6191
$ws.Range("M63").Value = -1647.6191
$ws.Range("H66").Value = 3019.8845
$ws.Range("I66").Value = 2333.6191
$ws.Range("K66").Value = 11668.0955
$ws.Range("M66").Value = -8236.095499999999
$ws.Range("H74").Value = 2852.5186
$ws.Range("I74").Value = 2353
$ws.Range("J74").Value = 4279.7144
$ws.Range("K74").Value = 2353
$ws.Range("L74").Value = 4279.7144
$ws.Range("M74").Value = -1479
$ws.Range("N74").Value = -6027.7144
$ws.Range("H77").Value = 2852.5186
$ws.Range("I77").Value = 2353
$ws.Range("J77").Value = 4279.7144
$ws.Range("K77").Value = 11765
$ws.Range("L77").Value = 21398.572
$ws.Range("M77").Value = -7397
$ws.Range("N77").Value = -30134.572
$ws.Range("H97").Value = 544.5
$ws.Range("I97").Value = 454.75
$ws.Range("J97").Value = 1083
$ws.Range("K97").Value = 454.75
$ws.Range("L97").Value = 1083
$ws.Range("M97").Value = 41.25
$ws.Range("N97").Value = -2075
$ws.Range("H110").Value = 33094.965
$ws.Range("I110").Value = 39520.668
$ws.Range("J110").Value = 2251.6
$ws.Range("K110").Value = 39520.668
$ws.Range("L110").Value = 2251.6
$ws.Range("M110").Value = -37475.668
$ws.Range("N110").Value = -6341.6
$ws.Range("H122").Value = 31499.678
$ws.Range("I122").Value = 2447.238
$ws.Range("J122").Value = 92509.8
$ws.Range("K122").Value = 7341.714
$ws.Range("L122").Value = 277529.4
$ws.Range("M122").Value = -4891.714
$ws.Range("N122").Value = -282429.4
$ws.Range("H123").Value = 84388.75
$ws.Range("J123").Value = 84388.75
$ws.Range("L123").Value = 84388.75
$ws.Range("N123").Value = -94188.75
$ws.Range("H132").Value = 10995.576
$ws.Range("I132").Value = 11491.896
$ws.Range("K132").Value = 34475.688
$ws.Range("M132").Value = -31945.688
$ws.Range("H135").Value = 96658.164
$ws.Range("J135").Value = 96658.164
$ws.Range("L135").Value = 96658.164
$ws.Range("N135").Value = -106798.164
$ws.Range("H136").Value = 4278.3584
$ws.Range("I136").Value = 2818.8684
$ws.Range("K136").Value = 8456.6052
$ws.Range("M136").Value = -5906.6052
$ws.Range("H137").Value = 55074.4
$ws.Range("I137").Value = 24921
$ws.Range("J137").Value = 59713.383
$ws.Range("K137").Value = 24921
$ws.Range("L137").Value = 59713.383
$ws.Range("M137").Value = -19821
$ws.Range("N137").Value = -69913.383
$ws.Range("H138").Value = 122428.25
$ws.Range("J138").Value = 122428.25
$ws.Range("L138").Value = 122428.25
$ws.Range("N138").Value = -132708.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 70830.16
$ws.Range("I86").Value = 1555.8889
$ws.Range("J86").Value = 226697.25
$ws.Range("K86").Value = 1555.8889
$ws.Range("L86").Value = 226697.25
$ws.Range("M86").Value = -432.8888999999999
$ws.Range("N86").Value = -228943.25
$ws.Range("H89").Value = 70830.16
$ws.Range("I89").Value = 1555.8889
$ws.Range("J89").Value = 226697.25
$ws.Range("K89").Value = 7779.4445
$ws.Range("L89").Value = 1133486.25
$ws.Range("M89").Value = -2163.4445
$ws.Range("N89").Value = -1144718.25
$ws.Range("H94").Value = 1162.862
$ws.Range("I94").Value = 848.8095
$ws.Range("J94").Value = 1987.25
$ws.Range("K94").Value = 848.8095
$ws.Range("L94").Value = 1987.25
$ws.Range("M94").Value = -397.8095
$ws.Range("N94").Value = -2889.25
$ws.Range("H103").Value = 14599.7
$ws.Range("J103").Value = 14599.7
$ws.Range("L103").Value = 14599.7
$ws.Range("N103").Value = -16943.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 199.2
$ws.Range("I7").Value = 147.6
$ws.Range("K7").Value = 147.6
$ws.Range("M7").Value = -34.59999999999999
$ws.Range("H16").Value = 1606.0476
$ws.Range("I16").Value = 639.4
$ws.Range("J16").Value = 4022.6667
$ws.Range("K16").Value = 639.4
$ws.Range("L16").Value = 4022.6667
$ws.Range("M16").Value = -352.4
$ws.Range("N16").Value = -4596.6667
$ws.Range("H31").Value = 2634.9387
$ws.Range("I31").Value = 2335
$ws.Range("K31").Value = 2335
$ws.Range("M31").Value = -2040
$ws.Range("H34").Value = 2634.9387
$ws.Range("I34").Value = 2335
$ws.Range("K34").Value = 2335
$ws.Range("M34").Value = -2133
$ws.Range("H99").Value = 5766.1924
$ws.Range("I99").Value = 4866.3335
$ws.Range("J99").Value = 7790.875
$ws.Range("K99").Value = 4866.3335
$ws.Range("L99").Value = 7790.875
$ws.Range("M99").Value = -3368.3335
$ws.Range("N99").Value = -10786.875
$ws.Range("H105").Value = 596.7727
$ws.Range("I105").Value = 615.7619
$ws.Range("J105").Value = 198
$ws.Range("K105").Value = 615.7619
$ws.Range("L105").Value = 198
$ws.Range("M105").Value = 1131.2381
$ws.Range("N105").Value = -3692
$ws.Range("H113").Value = 1606.0476
$ws.Range("I113").Value = 639.4
$ws.Range("J113").Value = 4022.6667
$ws.Range("K113").Value = 639.4
$ws.Range("L113").Value = 4022.6667
$ws.Range("M113").Value = 1530.6
$ws.Range("N113").Value = -8362.6667
$ws.Range("H126").Value = 5766.1924
$ws.Range("I126").Value = 4866.3335
$ws.Range("J126").Value = 7790.875
$ws.Range("K126").Value = 14599.0005
$ws.Range("L126").Value = 23372.625
$ws.Range("M126").Value = -12129.0005
$ws.Range("N126").Value = -28312.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 75.666664
$ws.Range("I38").Value = 51.23077
$ws.Range("J38").Value = 115.375
$ws.Range("K38").Value = 153.69231
$ws.Range("L38").Value = 346.125
$ws.Range("M38").Value = 193.30769
$ws.Range("N38").Value = -1040.125
$ws.Range("H51").Value = 2249.5
$ws.Range("I51").Value = 1500
$ws.Range("J51").Value = 2999
$ws.Range("K51").Value = 4500
$ws.Range("L51").Value = 8997
$ws.Range("M51").Value = -4040
$ws.Range("N51").Value = -9917
$ws.Range("H129").Value = 1252.091
$ws.Range("J129").Value = 1768.8
$ws.Range("L129").Value = 5306.4
$ws.Range("N129").Value = -15306.4
$ws.Range("H130").Value = 1120.6666
$ws.Range("I130").Value = 1120.6666
$ws.Range("K130").Value = 3361.9998
$ws.Range("M130").Value = 1658.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 4988.4287
$ws.Range("I21").Value = 4988.4287
$ws.Range("K21").Value = 4988.4287
$ws.Range("M21").Value = -4815.4287
$ws.Range("H30").Value = 4988.4287
$ws.Range("I30").Value = 4988.4287
$ws.Range("K30").Value = 4988.4287
$ws.Range("M30").Value = -4883.4287
$ws.Range("H70").Value = 5127.1665
$ws.Range("I70").Value = 5076.75
$ws.Range("K70").Value = 5076.75
$ws.Range("M70").Value = -4806.75
$ws.Range("H73").Value = 5127.1665
$ws.Range("I73").Value = 5076.75
$ws.Range("K73").Value = 5076.75
$ws.Range("M73").Value = -4140.75
$ws.Range("H102").Value = 30284.805
$ws.Range("I102").Value = 36303.82
$ws.Range("J102").Value = 10723
$ws.Range("K102").Value = 36303.82
$ws.Range("L102").Value = 10723
$ws.Range("M102").Value = -34681.82
$ws.Range("N102").Value = -13967
$ws.Range("H109").Value = 25498.9
$ws.Range("J109").Value = 25498.9
$ws.Range("L109").Value = 25498.9
$ws.Range("N109").Value = -27578.9
$ws.Range("H122").Value = 2695.4583
$ws.Range("I122").Value = 2143.4443
$ws.Range("K122").Value = 6430.3329
$ws.Range("M122").Value = -3980.3329
$ws.Range("H126").Value = 82087.3
$ws.Range("J126").Value = 4798.6665
$ws.Range("L126").Value = 14395.9995
$ws.Range("N126").Value = -19335.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4907.357
$ws.Range("I16").Value = 2870.3
$ws.Range("J16").Value = 10000
$ws.Range("K16").Value = 2870.3
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = -2700.3
$ws.Range("N16").Value = -10340
$ws.Range("H22").Value = 8223.200000000001
$ws.Range("I22").Value = 3633
$ws.Range("J22").Value = 8733.223
$ws.Range("K22").Value = 3633
$ws.Range("L22").Value = 8733.223
$ws.Range("M22").Value = -3338
$ws.Range("N22").Value = -9323.223
$ws.Range("H27").Value = 8223.200000000001
$ws.Range("I27").Value = 3633
$ws.Range("J27").Value = 8733.223
$ws.Range("K27").Value = 3633
$ws.Range("L27").Value = 8733.223
$ws.Range("M27").Value = -3526
$ws.Range("N27").Value = -8947.223
$ws.Range("H46").Value = 2388.4285
$ws.Range("I46").Value = 1054.75
$ws.Range("J46").Value = 4166.6665
$ws.Range("K46").Value = 1054.75
$ws.Range("L46").Value = 4166.6665
$ws.Range("M46").Value = -866.75
$ws.Range("N46").Value = -4542.6665
$ws.Range("H61").Value = 2702.5
$ws.Range("I61").Value = 2428.3125
$ws.Range("K61").Value = 2428.3125
$ws.Range("M61").Value = -2226.3125
$ws.Range("H68").Value = 3743.8
$ws.Range("I68").Value = 2496.4167
$ws.Range("J68").Value = 8733.333000000001
$ws.Range("K68").Value = 2496.4167
$ws.Range("L68").Value = 8733.333000000001
$ws.Range("M68").Value = -1747.4167
$ws.Range("N68").Value = -10231.333
$ws.Range("H71").Value = 3743.8
$ws.Range("I71").Value = 2496.4167
$ws.Range("J71").Value = 8733.333000000001
$ws.Range("K71").Value = 12482.0835
$ws.Range("L71").Value = 43666.665
$ws.Range("M71").Value = -8738.083500000001
$ws.Range("N71").Value = -51154.665
$ws.Range("H87").Value = 50000
$ws.Range("I87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("K87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("M87").Value = -48877
$ws.Range("N87").Value = -52246
$ws.Range("H90").Value = 50000
$ws.Range("I90").Value = 50000
$ws.Range("J90").Value = 50000
$ws.Range("K90").Value = 150000
$ws.Range("L90").Value = 150000
$ws.Range("M90").Value = -144384
$ws.Range("N90").Value = -161232
$ws.Range("H113").Value = 2702.5
$ws.Range("I113").Value = 2428.3125
$ws.Range("K113").Value = 2428.3125
$ws.Range("M113").Value = -258.3125
$ws.Range("H132").Value = 2908.3584
$ws.Range("I132").Value = 2659.1333
$ws.Range("K132").Value = 7977.3999
$ws.Range("M132").Value = -5447.3999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 139739.2
$ws.Range("I62").Value = 255636.75
$ws.Range("J62").Value = 7284.857
$ws.Range("K62").Value = 255636.75
$ws.Range("L62").Value = 7284.857
$ws.Range("M62").Value = -255012.75
$ws.Range("N62").Value = -8532.857
$ws.Range("H65").Value = 139739.2
$ws.Range("I65").Value = 255636.75
$ws.Range("J65").Value = 7284.857
$ws.Range("K65").Value = 1278183.75
$ws.Range("L65").Value = 36424.285
$ws.Range("M65").Value = -1275063.75
$ws.Range("N65").Value = -42664.285
$ws.Range("H88").Value = 85666
$ws.Range("J88").Value = 85666
$ws.Range("L88").Value = 85666
$ws.Range("N88").Value = -86478
$ws.Range("H91").Value = 85666
$ws.Range("J91").Value = 85666
$ws.Range("L91").Value = 85666
$ws.Range("N91").Value = -88474
$ws.Range("H96").Value = 3306.5715
$ws.Range("I96").Value = 3709.6
$ws.Range("J96").Value = 2299
$ws.Range("K96").Value = 3709.6
$ws.Range("L96").Value = 2299
$ws.Range("M96").Value = -2336.6
$ws.Range("N96").Value = -5045
$ws.Range("H100").Value = 9617270
$ws.Range("I100").Value = 12821027
$ws.Range("J100").Value = 5998
$ws.Range("K100").Value = 25642054
$ws.Range("L100").Value = 11996
$ws.Range("M100").Value = -25641513
$ws.Range("N100").Value = -13078
$ws.Range("H113").Value = 1330.325
$ws.Range("I113").Value = 1258
$ws.Range("K113").Value = 3774
$ws.Range("M113").Value = -1604
$ws.Range("H123").Value = 99999
$ws.Range("J123").Value = 99999
$ws.Range("L123").Value = 99999
$ws.Range("N123").Value = -109799
$ws.Range("H126").Value = 23862.21
$ws.Range("I126").Value = 31486.285
$ws.Range("K126").Value = 94458.855
$ws.Range("M126").Value = -91988.855
$ws.Range("H132").Value = 2780.116
$ws.Range("I132").Value = 1877.325
$ws.Range("J132").Value = 4025.3447
$ws.Range("K132").Value = 5631.975
$ws.Range("L132").Value = 12076.0341
$ws.Range("M132").Value = -3101.975
$ws.Range("N132").Value = -17136.0341
$ws.Range("H133").Value = 97499
$ws.Range("J133").Value = 97499
$ws.Range("L133").Value = 97499
$ws.Range("N133").Value = -107619
$ws.Range("H136").Value = 2918.7708
$ws.Range("I136").Value = 2689.9546
$ws.Range("K136").Value = 8069.8638
$ws.Range("M136").Value = -5519.8638

Write-Output "Updated market figures across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets."
